# Olyve.xlsx - "Others" sheet updates
# - Header E1: "Invalid Zip Code Characters" -> "Unsupported Zip Code"
# - Data row 2:
#     D2: 99387 (number)            -> "ABCDE" (text)
#     E2: "ABCDE" (text)            -> 11865 (number)
#     I2: (empty)                   -> "Please match the requested format"
#     K2: "...Reference: b0ba9c0add33" -> "...Reference:" (trailing hash removed)
# - Column I width widened (best-fit for the new, longer header text)
# - Selection moved to K2, with the view scrolled right so column C is
#   the left-most visible column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data row (set I2 before E1 so new shared strings are appended in the
#     same order as the authoritative edit) ---
$ws.Range("D2").Value = "ABCDE"
$ws.Range("E2").Value = 11865
$ws.Range("I2").Value = "Please match the requested format"

# --- Header row ---
$ws.Range("E1").Value = "Unsupported Zip Code"

$ws.Range("K2").Value = "Oops! There is a problem. Please enter a valid zip code. Reference:"

# --- Column width for column I (9th column) to fit the new header/value ---
$ws.Columns.Item(9).ColumnWidth = 32

# --- View state: scroll so column C is left-most, select K2 ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("K2").Select()
